$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column A (everything else shifts right by one).
$ws.Range("A1").EntireColumn.Insert()

# 2. Give the new header cell (A1) the same formatting as the other header
#    cells (bold font + fill), by copying the format from the neighboring
#    header cell (B1, which used to be A1).
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Fill in the new "RefID" column. The values are entered in this specific
#    order (RefID, LATFLD-26, LATFLD-27, LATFLD-25, LATFLD-28) to mirror how
#    they were originally typed in.
$ws.Range("A1").Value = "RefID"
$ws.Range("A2").Value = "LATFLD-26"
$ws.Range("A5").Value = "LATFLD-27"
$ws.Range("A3").Value = "LATFLD-25"
$ws.Range("A4").Value = "LATFLD-28"

# 4. Size the new column similarly to the other "best fit" columns.
$ws.Columns.Item(1).ColumnWidth = 9.17

# 5. Recreate the (hidden) AutoFilter defined name for the new used range.
$name = $ws.Names.Add("_xlnm._FilterDatabase", "=Transmittals_Close_Cancel!`$B`$1:`$Q`$5")
$name.Visible = $false

# 6. Reset the view back to the top-left/default selection.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A1").Select() | Out-Null
